# ACC Bid.xlsx — add "view history" and "edit profil" task sheets
# (commit: "task history + edit profil")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create the two new sheets at the end of the workbook. "edit profil" is
# created first (lower internal sheetId), "view history" second (higher
# sheetId), then "view history" is moved in front of "edit profil" so the
# final tab order is: ... task 23, view history, edit profil — matching
# sheetId=9 (view history) / sheetId=8 (edit profil) from the target.
#
# NOTE: this COM host re-resolves worksheet variable handles by *position*
# after a Move(), so handles captured before the move silently point at
# the wrong sheet afterwards. Always re-fetch sheets by name once the
# final tab order has been established, before touching their contents.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tmp1 = $wb.Worksheets.Add($null, $lastSheet)
$tmp1.Name = "edit profil"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$tmp2 = $wb.Worksheets.Add($null, $lastSheet2)
$tmp2.Name = "view history"

$tmp2.Move($tmp1)

# Re-fetch clean, final references by name (handles above are now stale).
$viewHistory = $wb.Worksheets.Item("view history")
$editProfil = $wb.Worksheets.Item("edit profil")

# ---------------------------------------------------------------------------
# "view history" sheet content
# ---------------------------------------------------------------------------
$viewHistory.Range("A1").Font.Bold = $true
$viewHistory.Range("A1").Value = "view"
$viewHistory.Range("B1").Font.Bold = $true
$viewHistory.Range("B1").Value = "expected"

$viewHistory.Range("A2").Value = "event sedang bejalan"
$viewHistory.Range("B2").Value = "pass"

$viewHistory.Range("A3").Value = "event selesai"
$viewHistory.Range("B3").Value = "pass"

$viewHistory.Columns.Item(1).ColumnWidth = 19.31
$null = $viewHistory.Range("B4").Select()
$viewHistory.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# "edit profil" sheet content
# ---------------------------------------------------------------------------
$editProfil.Range("A1").Font.Bold = $true
$editProfil.Range("A1").Value = "edit"

$editProfil.Range("B1").NumberFormat = "@"
$editProfil.Range("B1").Font.Bold = $true
$editProfil.Range("B1").Value = "varInput"

$editProfil.Range("C1").Font.Bold = $true
$editProfil.Range("C1").Value = "expected"

$editProfil.Range("A2").Value = "foto profil"
$editProfil.Range("B2").NumberFormat = "@"
$editProfil.Range("B2").Value = "IMG_20200128_092536.jpg"
$editProfil.Range("C2").Value = "pass"

$editProfil.Range("A3").Value = "no handphone"
$editProfil.Range("B3").NumberFormat = "@"
$editProfil.Range("B3").Value = "081325184829"
$editProfil.Range("C3").Value = "fail"

$editProfil.Range("A4").Value = "no handphone"
$editProfil.Range("B4").NumberFormat = "@"
$editProfil.Range("B4").Value = "085242869607"
$editProfil.Range("C4").Value = "pass"

$editProfil.Range("A5").Value = "pekerjaan"
$editProfil.Range("B5").NumberFormat = "@"
$editProfil.Range("B5").Value = "Wiraswasta"
$editProfil.Range("C5").Value = "pass"

$editProfil.Range("A6").Value = "alamat ktp"
$editProfil.Range("B6").NumberFormat = "@"
$editProfil.Range("B6").Value = "Yogyakarta"
$editProfil.Range("C6").Value = "pass"

$editProfil.Range("A7").Value = "tgl lahir"
$editProfil.Range("B7").NumberFormat = "@"
$editProfil.Range("B7").Value = "11/11/1998"
$editProfil.Range("C7").Value = "pass"

$editProfil.Range("A8").Value = "tgl lahir"
$editProfil.Range("B8").NumberFormat = "@"
$editProfil.Range("B8").Value = "12/01/2020"
$editProfil.Range("C8").Value = "fail"

$editProfil.Columns.Item(1).ColumnWidth = 13.17
$editProfil.Columns.Item(2).ColumnWidth = 23.45

$null = $editProfil.Range("B2").Select()
$editProfil.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Make "edit profil" the active tab (re-fetch by name once more — Select()
# calls above don't disturb handles, but re-fetching keeps this robust).
# ---------------------------------------------------------------------------
$finalEditProfil = $wb.Worksheets.Item("edit profil")
$finalEditProfil.Activate()
